$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column P: width + header "DataNascimento" ---
$ws.Columns.Item(16).ColumnWidth = 13.14
$ws.Range("O2").Copy($ws.Range("P2"))
$ws.Range("P2").Value = "DataNascimento"

# --- Extend the "Usuarios" title merge from I1:O1 to I1:P1 ---
$ws.Range("I1:O1").UnMerge()
$ws.Range("I1").Value = "Usuarios "
$ws.Range("I1").Copy($ws.Range("P1"))
$ws.Range("I1:P1").Merge()

# --- Fill P3:P7 with DataNascimento values (dates), matching I3:I7 styling ---
$ws.Range("I3").Copy($ws.Range("P3"))
$ws.Range("P3").Value = "02/02/2000"
$ws.Range("P3").NumberFormat = "mm-dd-yy"
$ws.Range("P3").Copy($ws.Range("P4:P7"))
$ws.Range("P4").Value = "02/02/2001"
$ws.Range("P5").Value = "02/02/2002"
$ws.Range("P6").Value = "02/02/2003"
$ws.Range("P7").Value = "02/02/2004"

# --- Update view: scroll so column C is leftmost, select M15 ---
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("M15").Select()
